$d = $word.ActiveDocument

function Set-IlvlZero($para) {
    # Adds <w:ilvl w:val="0"/> to a paragraph's numPr while leaving numId
    # untouched (it stays at 0, i.e. "not really in a list"). Mirrors the
    # ListLevelNumber + RemoveNumbers COM idiom Word itself uses for this.
    $para.Range.ListFormat.ListLevelNumber = 1
    $para.Range.ListFormat.RemoveNumbers()
}

# ---------------------------------------------------------------------
# 1) Three paragraphs whose <w:numPr><w:numId w:val="0"/></w:numPr> needs
#    an added <w:ilvl w:val="0"/>: "Second day(tuesday):", the blank line
#    right before "Third day (wednesday):", and "Third day (wednesday):"
#    itself.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Second day*tuesday*") {
        Set-IlvlZero $p
    }
}

$paras = @($d.Paragraphs)
for ($i = 0; $i -lt $paras.Count; $i++) {
    $t = $paras[$i].Range.Text
    if ($t -like "Third day*wednesday*") {
        Set-IlvlZero $paras[$i]
        # the blank paragraph immediately preceding it
        Set-IlvlZero $paras[$i - 1]
    }
}

# ---------------------------------------------------------------------
# 2) "Watching videos ... user security" + bookmark + "." -> merge into a
#    single run "Watching videos ... user security." and drop the
#    now-orphaned _GoBack bookmark from this paragraph (it is recreated
#    on the new "Fourth day (Thursday)" paragraph below).
# ---------------------------------------------------------------------
$watching = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Watching videos*") {
        $watching = $p
    }
}

$paraEnd = $watching.Range.End - 1
$tailRange = $d.Range($paraEnd - 3, $paraEnd)
$tailRange.Delete()
$reinsertAt = $d.Range($watching.Range.Start, $watching.Range.End - 1)
$reinsertAt.InsertAfter("ty.")

# ---------------------------------------------------------------------
# 3) Insert a new blank paragraph after "Watching videos ..." and before
#    the old trailing blank paragraph, matching the target pPr.
# ---------------------------------------------------------------------
$trailing = $d.Paragraphs.Last
$insAt = $d.Range($trailing.Range.Start, $trailing.Range.Start)
$blankFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr><w:ind w:leftChars="0"/><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$insAt.InsertXML($blankFrag)

# ---------------------------------------------------------------------
# 4) Turn the (still-present, now-last) old trailing blank paragraph into
#    "Fourth day (Thursday)" carrying the _GoBack bookmark, and give it
#    <w:ilvl w:val="0"/> too.
# ---------------------------------------------------------------------
$fourth = $d.Paragraphs.Last
$fourthAt = $d.Range($fourth.Range.Start, $fourth.Range.Start)
$fourthFrag = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US"/></w:rPr><w:t>Fourth day (Thursday)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$fourthAt.InsertXML($fourthFrag)
Set-IlvlZero $fourth
